$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Fitness column (C) for the first generations of Run 28.
# Rows 2-12 correspond to Generation 0..10 -> Fitness 7312
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 3).Value = 7312
}

# Rows 13-14 correspond to Generation 11..12 -> Fitness 7295
for ($r = 13; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = 7295
}
